$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header (row1) / data (row2) cells in the exact order the author entered them,
# which determines the order they land in the shared string table.
$ws.Range("AL1").Value = "Modifier name"
$ws.Range("AL2").Value = "Test001, Test002"

$ws.Range("AM1").Value = "Modifier Type"
$ws.Range("AN1").Value = "Modifier Value"
$ws.Range("AN2").Value = "Test prodouct client voucher, Test product key - fulfilled ebook - component"
$ws.Range("AM2").Value = "Dropdown, Dropdown"

$ws.Range("AO1").Value = "Required"
$ws.Range("AO2").Value = "Yes, No"

# Update existing product name / sku in row 2 (these come last in the shared strings table)
$ws.Range("A2").Value = "Comptia Product 54 Voucher"
$ws.Range("B2").Value = "CP-1054"

# Copy header style from existing header cell (AK1) to new header cells
$ws.Range("AK1").Copy()
$ws.Range("AL1:AO1").PasteSpecial(-4122)

# Apply wrap text style to new data cells
$ws.Range("AL2:AO2").WrapText = $true

# Adjust row height for row2
$ws.Rows("2").RowHeight = 158.4

# Adjust view (scroll so row 2 is the top visible row, then select F9)
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F9").Select()
